# Add season-record columns (Wins / Losses / Ties) to the player table.
# New columns AD:AF are appended after the existing AC ("Unnamed: 28") column,
# extending the used range from A1:AC48 to A1:AF48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels, styled like the rest of the header ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous

# --- Data rows (2-48): every player gets the team's season record ---
for ($row = 2; $row -le 48; $row++) {
    $ws.Cells.Item($row, 30).Value = 75   # AD -> Wins
    $ws.Cells.Item($row, 31).Value = 87   # AE -> Losses
    $ws.Cells.Item($row, 32).Value = 0    # AF -> Ties
}
